# Commit: "lls time and poe2 prov"
# Adds a new LLS\Provider test case ("L_Prov_Remitter_To") row to both the
# summary sheet (Sheet1) and the detail sheet (Sheet2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1 - summary
$ws2 = $wb.Worksheets.Item(2)   # Sheet2 - detail
$ws3 = $wb.Worksheets.Item(3)   # Sheet3 - untouched, kept active-state as-is

# ---------------------------------------------------------------------------
# Sheet2: append the new detail row right after the current last row (35)
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("D36").Value = "LLS\Provider"
$ws2.Range("E36").Value = "L_Prov_Remitter_To"
$ws2.Range("D36:E36").Select()

# ---------------------------------------------------------------------------
# Sheet1: append the matching summary row right after the current last row (2)
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A3").Value = "LLS\Provider"
$ws1.Range("B3").Value = "L_Prov_Remitter_To"
$ws1.Range("A4:XFD33").Select()

Write-Host "Added L_Prov_Remitter_To rows to Sheet1 and Sheet2"
